$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are textual in the source sheet (e.g. "30.426.81",
# "1.000"). Force text format before assignment so Excel COM does not
# reinterpret them as numbers, then restore the default "Normal" style so no
# stray formatting diff is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.426.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4804"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2788"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.863.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07441"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.074"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6384"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.396.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").Value = "  -2.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "232.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007462"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.118.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.124"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.078"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.303"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.896"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1045"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.379"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.254"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.965"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04964"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.167"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7390"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01935"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.631"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9136"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.032"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("E44").Value = "  -1.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.575"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.154"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1223"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.844"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.413"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "

